$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert a new "Meta description" paragraph right after the title
# (Heading1) paragraph. We first split the title paragraph into two
# paragraphs via Find/Replace (scoped to the title paragraph's own Range so
# we don't touch the later duplicate of the same text), which leaves a
# brand-new, empty paragraph right after the title. We then use InsertXML on
# that paragraph's Range to stamp the exact run structure we need (a leading
# empty run, a bold "Meta description" run, and a plain run with the rest of
# the sentence) - this also clears the inherited Heading1 paragraph
# properties, since InsertXML replaces the whole range (including the
# paragraph mark) with our own markup.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.First
$titleText = "Play Doom of Egypt for Free - Review 2021 | AP"
$titleRng = $titlePara.Range
$titleRng.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, ($titleText + "^p"), 2)

$metaPara = $d.Paragraphs(2)
$metaRng = $metaPara.Range

$metaBold = "Meta description"
$metaRest = ": Read our review of Doom of Egypt slot game and play for free. Discover what we like and don't like about the game's visuals, payouts, volatility, and theme."

$metaXml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?>" +
    "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
    "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:body><w:p>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>$metaBold</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>$metaRest</w:t></w:r>" +
    "</w:p></w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

$metaRng.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# Step 2: Remove the duplicate bold title paragraph that was sitting just
# before the final (italic) meta-description paragraph at the end of the
# document.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($i -ne 1 -and $text -eq $titleText) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Step 3: Replace the text of the final (italic) paragraph - the old SEO
# meta-description placeholder - with the new image-generation prompt text.
# We target only the run's own text range (not the whole paragraph range,
# which also holds the leading empty run) and assign Range.Text directly so
# that no smart-quote autocorrection is applied to the straight quotes in
# the new text.
# ---------------------------------------------------------------------------

$oldPromptText = "Read our review of Doom of Egypt slot game and play for free. Discover what we like and don't like about the game's visuals, payouts, volatility, and theme."
$newPromptText = 'Prompt: Create a cartoon-style feature image that represents the game "Doom of Egypt" and features a happy Maya warrior with glasses. The character should be dressed in traditional Maya clothing and holding a golden scarab or an Anubis symbol. The background should be set inside a pyramid, with a misty blue atmosphere and hieroglyphics adorning the walls. The overall design should be eye-catching and engaging, highlighting the game''s macabre yet captivating theme.'

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRng = $lastPara.Range
$promptRng = $d.Range($lastRng.Start, $lastRng.Start + $oldPromptText.Length)
$promptRng.Text = $newPromptText

Write-Output "done"
